# Deploy IG build 2e04aec7a6e183fa6df6ea6f057cb991a5ed746d
#
# Updates the "Metadata" sheet of the FHIR IG ValueSet export with the
# refreshed build metadata: Experimental flag, Date, and Description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Property"/"Value" rows on the Metadata sheet are plain text cells
# (style index 2 — border + top-aligned wrap text, General format). Excel's
# autoformatting would otherwise turn "false" into a Boolean and the ISO
# date into a date serial, so each value is entered with a leading
# apostrophe to keep it as literal text, and then the original cell
# formatting (copied from an already-correct text cell in the same column)
# is reapplied so the style index is not disturbed.

# Experimental: no value -> "false"
$ws.Range("B7").Value = "'false"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date: regenerated build date (kept as literal text, not a date serial)
$ws.Range("B8").Value = "'2023-10-26"
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Description: no value -> the value set's descriptive text
$ws.Range("B11").Value = "Kodeverket er for å angi dokumenttype når man skal dele et dokument i en XDS-løsning. Kodeverket er også harmonert med eksisterende kodeverk for dokumenttyper og journalstruktur i den grad dette har vært mulig. Dette gjelder primært følgende kodeverk: - 9066 Kategori journalinformasjon - 9601 Emneorienterte sakstyper (Piene inndelingen)"

$excel.CutCopyMode = $false
